$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: C2 and E2 values were removed (cells cleared / emptied)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: C3 removed, E3 updated
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = -1.563633964192079

# Row 4: C4 updated
$ws.Range("C4").Value = -4.774178217057756

# Row 6: E6 updated
$ws.Range("E6").Value = 2.114249845651872

# Row 7: C7 updated
$ws.Range("C7").Value = 1.239479831392831

# Row 8: C8 updated
$ws.Range("C8").Value = 0.2379616621360992

# Row 10: C10, E10 updated
$ws.Range("C10").Value = 1.470039379455734
$ws.Range("E10").Value = 1.339087911421122

# Row 11: C11, E11 updated
$ws.Range("C11").Value = 1.638797242243228
$ws.Range("E11").Value = 1.006353890555212

# Row 13: E13 updated
$ws.Range("E13").Value = -0.301339632123987

# Row 14: E14 updated
$ws.Range("E14").Value = 0.2691345740890139

# Row 15: E15 updated
$ws.Range("E15").Value = 23.52713729381606

# Row 16: C16, E16 updated
$ws.Range("C16").Value = 1.099928004397577
$ws.Range("E16").Value = 6.182044950645027

# Row 17: C17 updated
$ws.Range("C17").Value = 2.310042359896247

# Row 18: E18 updated
$ws.Range("E18").Value = -0.3513551123189074

# Row 19: C19, E19 updated
$ws.Range("C19").Value = -0.3101476031197037
$ws.Range("E19").Value = 0.2561130241983456
